# [ ADD SUPPORT OF NO TOKEN OF DIAGNOSTIC SIGNALS AND DIAGNOSTIC FRAMES ]
#
# 1) "Frame and Signal Attributes": the "Signal representation" column (J)
#    no longer carries the per-signal encoding token strings - those rows
#    now just show "/" (the workbook's usual "not applicable" placeholder),
#    matching the format already used by every other "no value" cell.
# 2) "Diagnostic_Frames": the per-signal breakdown rows for the diagnostic
#    frames (MasterReq/SlaveResp and their byte-by-byte signals) are removed
#    entirely - only the header row is kept.

$wb = $excel.ActiveWorkbook

# --- 1) Frame and Signal Attributes: clear "Signal representation" tokens ---
$ws3 = $wb.Worksheets.Item("Frame and Signal Attributes")

$sigRepr = $ws3.Range("J2:J5")
$sigRepr.HorizontalAlignment = -4108   # xlCenter - match the rest of the row
$sigRepr.WrapText = $false
$sigRepr.Value = "/"

# --- 2) Diagnostic_Frames: drop the diagnostic signal rows, keep the header ---
$ws4 = $wb.Worksheets.Item("Diagnostic_Frames")
$ws4.Range("A2:A17").EntireRow.Delete()
